# Updated symbol list on Fri Feb 17 08:58:57 UTC 2023 with GitHub Actions
#
# The Price (D) / Volume(1h) (E) columns hold plain text that merely looks
# numeric ("309.50", "-3.69%", "2,022.23%", ...). Writing such a string to
# a .Value directly would let Excel auto-coerce it into a real number (and
# drop formatting like trailing zeros or thousands separators), so every
# D/E write below is prefixed with a leading apostrophe, exactly like a
# user typing a quote-prefixed literal in the Excel UI, to force it to
# stay text. Coin name (B) / link (C) columns are plain non-numeric text
# and need no such trick.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Cells.Item(2, 4).Value = "'310.39"
$ws.Cells.Item(2, 5).Value = "'-3.44%"

# Row 3 - OKB
$ws.Cells.Item(3, 4).Value = "'49.11"
$ws.Cells.Item(3, 5).Value = "'1.19%"

# Row 4 - HuobiToken
$ws.Cells.Item(4, 4).Value = "'5.153"
$ws.Cells.Item(4, 5).Value = "'-1.60%"

# Row 5 - Cronos
$ws.Cells.Item(5, 4).Value = "'0.07780"
$ws.Cells.Item(5, 5).Value = "'-3.97%"

# Row 6 - GateToken
$ws.Cells.Item(6, 4).Value = "'4.524"
$ws.Cells.Item(6, 5).Value = "'-1.64%"

# Row 7 - MXToken
$ws.Cells.Item(7, 4).Value = "'1.372"

# Row 8 - FTXToken
$ws.Cells.Item(8, 4).Value = "'1.568"
$ws.Cells.Item(8, 5).Value = "'-4.54%"

# Row 9 - LiechtensteinCryptoassetsExchange
$ws.Cells.Item(9, 4).Value = "'0.1223"
$ws.Cells.Item(9, 5).Value = "'-6.71%"

# Row 10 - WazirX
$ws.Cells.Item(10, 5).Value = "'0.33%"

# Row 11 - BitrueCoin
$ws.Cells.Item(11, 4).Value = "'0.04706"
$ws.Cells.Item(11, 5).Value = "'2.55%"

# Row 12 - MandalaExchangeToken
$ws.Cells.Item(12, 4).Value = "'0.09391"
$ws.Cells.Item(12, 5).Value = "'-1.70%"

# Row 13 - BitMartToken
$ws.Cells.Item(13, 4).Value = "'0.1043"
$ws.Cells.Item(13, 5).Value = "'-0.61%"

# Row 14 - BitForexToken
$ws.Cells.Item(14, 4).Value = "'0.001263"
$ws.Cells.Item(14, 5).Value = "'-5.79%"

# Row 15 - CoinExToken
$ws.Cells.Item(15, 4).Value = "'0.04178"
$ws.Cells.Item(15, 5).Value = "'-2.67%"

# Row 16 - TigerCash
$ws.Cells.Item(16, 4).Value = "'0.005823"
$ws.Cells.Item(16, 5).Value = "'-2.03%"

# Row 17 - UpBots (Price unchanged, only Volume(1h) moves)
$ws.Cells.Item(17, 5).Value = "'2,021.23%"

# Rows 18-25 shift down by one (a new HotbitToken entry is inserted at the
# top of the block), so every row's Coin/Link/Price/Volume gets replaced
# with what used to be the row above it (with its own Price/Volume refresh).

# Row 18 - now HotbitToken (was LEO)
$ws.Cells.Item(18, 2).Value = "HotbitToken"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(18, 4).Value = "'0.003922"
$ws.Cells.Item(18, 5).Value = "'-7.78%"

# Row 19 - now LEO (was BTSEToken)
$ws.Cells.Item(19, 2).Value = "LEO"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(19, 4).Value = "'3.330"
$ws.Cells.Item(19, 5).Value = "'-0.44%"

# Row 20 - now BTSEToken (was BitpandaEcosystemToken)
$ws.Cells.Item(20, 2).Value = "BTSEToken"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(20, 4).Value = "'2.431"
$ws.Cells.Item(20, 5).Value = "'-0.25%"

# Row 21 - now BitpandaEcosystemToken (was MCDex)
$ws.Cells.Item(21, 2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Cells.Item(21, 4).Value = "'0.3391"
$ws.Cells.Item(21, 5).Value = "'-0.38%"

# Row 22 - now MCDex (was ProBitToken)
$ws.Cells.Item(22, 2).Value = "MCDex"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(22, 4).Value = "'8.013"
$ws.Cells.Item(22, 5).Value = "'-2.48%"

# Row 23 - now ProBitToken (was ZBToken)
$ws.Cells.Item(23, 2).Value = "ProBitToken"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Cells.Item(23, 4).Value = "'0.1348"
$ws.Cells.Item(23, 5).Value = "'-4.41%"

# Row 24 - now ZBToken (was BitKan)
$ws.Cells.Item(24, 2).Value = "ZBToken"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Cells.Item(24, 4).Value = "'0.3039"
$ws.Cells.Item(24, 5).Value = "'-3.48%"

# Row 25 - now BitKan (was HotbitToken)
$ws.Cells.Item(25, 2).Value = "BitKan"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(25, 4).Value = "'0.001275"
$ws.Cells.Item(25, 5).Value = "'-2.25%"

# Row 26 - NitroEx (Coin/Link unchanged, only Price/Volume refresh)
$ws.Cells.Item(26, 4).Value = "'0.0001351"
$ws.Cells.Item(26, 5).Value = "'0.08%"

# Row 38 - One
$ws.Cells.Item(38, 4).Value = "'0.02586"
$ws.Cells.Item(38, 5).Value = "'-3.07%"

# Row 39 - IDEX
$ws.Cells.Item(39, 4).Value = "'0.06007"
$ws.Cells.Item(39, 5).Value = "'6.76%"

# Row 40 - CEJI
$ws.Cells.Item(40, 4).Value = "'0.01103"
$ws.Cells.Item(40, 5).Value = "'75.13%"

# Row 41 - KickToken
$ws.Cells.Item(41, 4).Value = "'0.007931"
$ws.Cells.Item(41, 5).Value = "'3.20%"

# Row 42 - BKEXToken
$ws.Cells.Item(42, 4).Value = "'0.1418"
$ws.Cells.Item(42, 5).Value = "'-1.44%"

# Row 43 - Dexo
$ws.Cells.Item(43, 4).Value = "'0.008432"
$ws.Cells.Item(43, 5).Value = "'9.80%"

# Row 44 - LocalTraders
$ws.Cells.Item(44, 4).Value = "'0.008333"
$ws.Cells.Item(44, 5).Value = "'2.86%"

# Row 45 - PooCoin
$ws.Cells.Item(45, 4).Value = "'0.3131"
$ws.Cells.Item(45, 5).Value = "'-1.97%"

# Row 46 - CoinLion
$ws.Cells.Item(46, 4).Value = "'0.00007663"
$ws.Cells.Item(46, 5).Value = "'9.41%"

# Row 47 - Kangarootoken
$ws.Cells.Item(47, 4).Value = "'0.00000000752"
$ws.Cells.Item(47, 5).Value = "'0.24%"

# Row 48 - BOLO (Price unchanged, only Volume(1h) moves)
$ws.Cells.Item(48, 5).Value = "'-0.66%"

# Row 49 - CoinbaseStockToken
$ws.Cells.Item(49, 4).Value = "'0.002626"
$ws.Cells.Item(49, 5).Value = "'-34.37%"

# Row 50 - CryptobidCoin
$ws.Cells.Item(50, 4).Value = "'0.00002105"
$ws.Cells.Item(50, 5).Value = "'0.24%"

# Row 51 - SpecialPowerGold
$ws.Cells.Item(51, 4).Value = "'0.0002005"
$ws.Cells.Item(51, 5).Value = "'0.24%"
